$d = $word.ActiveDocument

function New-PkgXml([string]$bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Find the paragraphs we need to touch by their current text, so the script
# is resilient to exact paragraph-index assumptions.
$pRaksha = $null
$pHiRaksha = $null
$pSentAck = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Raksha") { $pRaksha = $p }
    elseif ($t -eq "Hi Raksha,") { $pHiRaksha = $p }
    elseif ($t -eq "Sent you Ack file") { $pSentAck = $p }
}

# 1) Signature "Raksha" paragraph: wrap the run with spell-check proof marks.
if ($pRaksha -ne $null) {
    $xmlRaksha = New-PkgXml '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Raksha</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
    $pRaksha.Range.InsertXML($xmlRaksha)
}

# 2) "Hi Raksha," paragraph: split into "Hi " / Raksha (proofed) / ",".
#    Also drop the _GoBack bookmark that used to live here (it is recreated
#    below on the new paragraph), doing this BEFORE adding the new bookmark
#    keeps the bookmark id at 0 instead of being bumped to avoid a clash.
if ($pHiRaksha -ne $null) {
    $xmlHi = New-PkgXml '<w:p><w:r><w:t xml:space="preserve">Hi </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Raksha</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r></w:p>'
    $pHiRaksha.Range.InsertXML($xmlHi)
}

# 3) "Sent you Ack file" paragraph: split into runs, add " on github", and
#    move the _GoBack bookmark to the end of this paragraph.
if ($pSentAck -ne $null) {
    $xmlSent = New-PkgXml '<w:p><w:r><w:t xml:space="preserve">Sent you </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file</w:t></w:r><w:r><w:t xml:space="preserve"> on github</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
    $pSentAck.Range.InsertXML($xmlSent)
}
